$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 53 holding what used to be row 52's data (shift it down),
# then overwrite row 52 with the new weekly price entry.

# 1) Copy the current row 52 values into the new row 53, cell by cell.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(53, $col).Value2 = $ws.Cells.Item(52, $col).Value2
}
$ws.Cells.Item(53, 4).NumberFormat = $ws.Cells.Item(52, 4).NumberFormat

# 2) Update row 52 with the new data (new date + new K:R values).
$ws.Cells.Item(52, 4).Value2 = 44568

$ws.Cells.Item(52, 11).Value2 = 40000
$ws.Cells.Item(52, 12).Value2 = 40000
$ws.Cells.Item(52, 13).Value2 = 40000
$ws.Cells.Item(52, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(52, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value2 = 1600
$ws.Cells.Item(52, 17).Value2 = 25
$ws.Cells.Item(52, 18).Value2 = "Hortaliza"

$wb.Save()
